# Update the P2P history report (Binance) with the latest export:
#  - refreshed BUY-side rows 2-10 (columns A:E) and SELL-side rows 2-14 (columns H:M)
#  - rows 11-12 no longer have matching BUY-side entries, so that data is cleared
#  - appended new SELL-side rows 15-40 from the latest pull
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove BUY-side (A:E) data for rows 11-12; only SELL-side (H:L) remains
$ws.Range("A11:E12").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = 'COMPLETED'
$ws.Cells.Item(2, 2).Value = 1141
$ws.Cells.Item(2, 3).Value = 'BUY'
$ws.Cells.Item(2, 4).Value = '2024-06-10 23:02:15'
$ws.Cells.Item(2, 8).Value = 'COMPLETED'
$ws.Cells.Item(2, 9).Value = 2044
$ws.Cells.Item(2, 10).Value = 'SELL'
$ws.Cells.Item(2, 11).Value = '2024-06-10 23:42:09'
$ws.Cells.Item(2, 13).Value = 56872.89999999999

# Row 3
$ws.Cells.Item(3, 1).Value = 'COMPLETED'
$ws.Cells.Item(3, 2).Value = 1094
$ws.Cells.Item(3, 3).Value = 'BUY'
$ws.Cells.Item(3, 4).Value = '2024-06-10 22:52:38'
$ws.Cells.Item(3, 8).Value = 'CANCELLED'
$ws.Cells.Item(3, 9).Value = 1200
$ws.Cells.Item(3, 10).Value = 'SELL'
$ws.Cells.Item(3, 11).Value = '2024-06-10 23:38:53'

# Row 4
$ws.Cells.Item(4, 1).Value = 'COMPLETED'
$ws.Cells.Item(4, 2).Value = 1595
$ws.Cells.Item(4, 3).Value = 'BUY'
$ws.Cells.Item(4, 4).Value = '2024-06-10 22:37:07'
$ws.Cells.Item(4, 8).Value = 'COMPLETED'
$ws.Cells.Item(4, 9).Value = 2023
$ws.Cells.Item(4, 10).Value = 'SELL'
$ws.Cells.Item(4, 11).Value = '2024-06-10 23:19:51'

# Row 5
$ws.Cells.Item(5, 1).Value = 'COMPLETED'
$ws.Cells.Item(5, 2).Value = 2953
$ws.Cells.Item(5, 3).Value = 'BUY'
$ws.Cells.Item(5, 4).Value = '2024-06-10 14:40:03'
$ws.Cells.Item(5, 8).Value = 'CANCELLED_BY_SYSTEM'
$ws.Cells.Item(5, 9).Value = 1200
$ws.Cells.Item(5, 10).Value = 'SELL'
$ws.Cells.Item(5, 11).Value = '2024-06-10 23:17:06'

# Row 6
$ws.Cells.Item(6, 1).Value = 'COMPLETED'
$ws.Cells.Item(6, 2).Value = 3561
$ws.Cells.Item(6, 3).Value = 'BUY'
$ws.Cells.Item(6, 4).Value = '2024-06-10 14:18:36'
$ws.Cells.Item(6, 8).Value = 'COMPLETED'
$ws.Cells.Item(6, 9).Value = 2000
$ws.Cells.Item(6, 10).Value = 'SELL'
$ws.Cells.Item(6, 11).Value = '2024-06-10 22:28:12'

# Row 7
$ws.Cells.Item(7, 1).Value = 'COMPLETED'
$ws.Cells.Item(7, 2).Value = 1201
$ws.Cells.Item(7, 3).Value = 'BUY'
$ws.Cells.Item(7, 4).Value = '2024-06-10 14:11:33'
$ws.Cells.Item(7, 8).Value = 'COMPLETED'
$ws.Cells.Item(7, 9).Value = 3000
$ws.Cells.Item(7, 10).Value = 'SELL'
$ws.Cells.Item(7, 11).Value = '2024-06-10 22:04:38'

# Row 8
$ws.Cells.Item(8, 1).Value = 'COMPLETED'
$ws.Cells.Item(8, 2).Value = 2600
$ws.Cells.Item(8, 3).Value = 'BUY'
$ws.Cells.Item(8, 4).Value = '2024-06-10 14:08:44'
$ws.Cells.Item(8, 8).Value = 'COMPLETED'
$ws.Cells.Item(8, 9).Value = 1500
$ws.Cells.Item(8, 10).Value = 'SELL'
$ws.Cells.Item(8, 11).Value = '2024-06-10 22:00:42'

# Row 9
$ws.Cells.Item(9, 1).Value = 'COMPLETED'
$ws.Cells.Item(9, 2).Value = 2122
$ws.Cells.Item(9, 3).Value = 'BUY'
$ws.Cells.Item(9, 4).Value = '2024-06-10 14:04:17'
$ws.Cells.Item(9, 8).Value = 'COMPLETED'
$ws.Cells.Item(9, 9).Value = 2510
$ws.Cells.Item(9, 10).Value = 'SELL'
$ws.Cells.Item(9, 11).Value = '2024-06-10 21:50:16'

# Row 10
$ws.Cells.Item(10, 1).Value = 'COMPLETED'
$ws.Cells.Item(10, 2).Value = 1343
$ws.Cells.Item(10, 3).Value = 'BUY'
$ws.Cells.Item(10, 4).Value = '2024-06-10 13:57:07'
$ws.Cells.Item(10, 8).Value = 'CANCELLED'
$ws.Cells.Item(10, 9).Value = 2500
$ws.Cells.Item(10, 10).Value = 'SELL'
$ws.Cells.Item(10, 11).Value = '2024-06-10 21:32:31'

# Row 11
$ws.Cells.Item(11, 8).Value = 'CANCELLED'
$ws.Cells.Item(11, 9).Value = 2000
$ws.Cells.Item(11, 10).Value = 'SELL'
$ws.Cells.Item(11, 11).Value = '2024-06-10 21:23:48'

# Row 12
$ws.Cells.Item(12, 8).Value = 'COMPLETED'
$ws.Cells.Item(12, 9).Value = 2424
$ws.Cells.Item(12, 10).Value = 'SELL'
$ws.Cells.Item(12, 11).Value = '2024-06-10 21:18:19'

# Row 13
$ws.Cells.Item(13, 8).Value = 'COMPLETED'
$ws.Cells.Item(13, 9).Value = 1200
$ws.Cells.Item(13, 10).Value = 'SELL'
$ws.Cells.Item(13, 11).Value = '2024-06-10 21:13:11'

# Row 14
$ws.Cells.Item(14, 8).Value = 'COMPLETED'
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 10).Value = 'SELL'
$ws.Cells.Item(14, 11).Value = '2024-06-10 21:12:05'

# Row 15
$ws.Cells.Item(15, 8).Value = 'COMPLETED'
$ws.Cells.Item(15, 9).Value = 1000
$ws.Cells.Item(15, 10).Value = 'SELL'
$ws.Cells.Item(15, 11).Value = '2024-06-10 21:11:46'

# Row 16
$ws.Cells.Item(16, 8).Value = 'CANCELLED_BY_SYSTEM'
$ws.Cells.Item(16, 9).Value = 4000
$ws.Cells.Item(16, 10).Value = 'SELL'
$ws.Cells.Item(16, 11).Value = '2024-06-10 21:09:32'

# Row 17
$ws.Cells.Item(17, 8).Value = 'CANCELLED'
$ws.Cells.Item(17, 9).Value = 1400
$ws.Cells.Item(17, 10).Value = 'SELL'
$ws.Cells.Item(17, 11).Value = '2024-06-10 21:02:49'

# Row 18
$ws.Cells.Item(18, 8).Value = 'COMPLETED'
$ws.Cells.Item(18, 9).Value = 1900
$ws.Cells.Item(18, 10).Value = 'SELL'
$ws.Cells.Item(18, 11).Value = '2024-06-10 20:59:25'

# Row 19
$ws.Cells.Item(19, 8).Value = 'CANCELLED_BY_SYSTEM'
$ws.Cells.Item(19, 9).Value = 1900
$ws.Cells.Item(19, 10).Value = 'SELL'
$ws.Cells.Item(19, 11).Value = '2024-06-10 20:51:19'

# Row 20
$ws.Cells.Item(20, 8).Value = 'CANCELLED'
$ws.Cells.Item(20, 9).Value = 3600
$ws.Cells.Item(20, 10).Value = 'SELL'
$ws.Cells.Item(20, 11).Value = '2024-06-10 20:37:12'

# Row 21
$ws.Cells.Item(21, 8).Value = 'COMPLETED'
$ws.Cells.Item(21, 9).Value = 1285
$ws.Cells.Item(21, 10).Value = 'SELL'
$ws.Cells.Item(21, 11).Value = '2024-06-10 19:37:29'

# Row 22
$ws.Cells.Item(22, 8).Value = 'CANCELLED'
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 'SELL'
$ws.Cells.Item(22, 11).Value = '2024-06-10 18:52:45'

# Row 23
$ws.Cells.Item(23, 8).Value = 'CANCELLED_BY_SYSTEM'
$ws.Cells.Item(23, 9).Value = 4000
$ws.Cells.Item(23, 10).Value = 'SELL'
$ws.Cells.Item(23, 11).Value = '2024-06-10 18:46:44'

# Row 24
$ws.Cells.Item(24, 8).Value = 'CANCELLED'
$ws.Cells.Item(24, 9).Value = 4000
$ws.Cells.Item(24, 10).Value = 'SELL'
$ws.Cells.Item(24, 11).Value = '2024-06-10 18:34:34'

# Row 25
$ws.Cells.Item(25, 8).Value = 'COMPLETED'
$ws.Cells.Item(25, 9).Value = 4000
$ws.Cells.Item(25, 10).Value = 'SELL'
$ws.Cells.Item(25, 11).Value = '2024-06-10 16:12:40'

# Row 26
$ws.Cells.Item(26, 8).Value = 'CANCELLED'
$ws.Cells.Item(26, 9).Value = 1700
$ws.Cells.Item(26, 10).Value = 'SELL'
$ws.Cells.Item(26, 11).Value = '2024-06-10 16:01:28'

# Row 27
$ws.Cells.Item(27, 8).Value = 'COMPLETED'
$ws.Cells.Item(27, 9).Value = 2000
$ws.Cells.Item(27, 10).Value = 'SELL'
$ws.Cells.Item(27, 11).Value = '2024-06-10 15:28:06'

# Row 28
$ws.Cells.Item(28, 8).Value = 'COMPLETED'
$ws.Cells.Item(28, 9).Value = 1000
$ws.Cells.Item(28, 10).Value = 'SELL'
$ws.Cells.Item(28, 11).Value = '2024-06-10 14:59:55'

# Row 29
$ws.Cells.Item(29, 8).Value = 'COMPLETED'
$ws.Cells.Item(29, 9).Value = 1000
$ws.Cells.Item(29, 10).Value = 'SELL'
$ws.Cells.Item(29, 11).Value = '2024-06-10 14:48:36'

# Row 30
$ws.Cells.Item(30, 8).Value = 'CANCELLED'
$ws.Cells.Item(30, 9).Value = 1000
$ws.Cells.Item(30, 10).Value = 'SELL'
$ws.Cells.Item(30, 11).Value = '2024-06-10 14:45:30'

# Row 31
$ws.Cells.Item(31, 8).Value = 'COMPLETED'
$ws.Cells.Item(31, 9).Value = 1000
$ws.Cells.Item(31, 10).Value = 'SELL'
$ws.Cells.Item(31, 11).Value = '2024-06-10 14:33:36'

# Row 32
$ws.Cells.Item(32, 8).Value = 'COMPLETED'
$ws.Cells.Item(32, 9).Value = 1000
$ws.Cells.Item(32, 10).Value = 'SELL'
$ws.Cells.Item(32, 11).Value = '2024-06-10 14:33:02'

# Row 33
$ws.Cells.Item(33, 8).Value = 'COMPLETED'
$ws.Cells.Item(33, 9).Value = 1600
$ws.Cells.Item(33, 10).Value = 'SELL'
$ws.Cells.Item(33, 11).Value = '2024-06-10 14:28:05'

# Row 34
$ws.Cells.Item(34, 8).Value = 'COMPLETED'
$ws.Cells.Item(34, 9).Value = 1000
$ws.Cells.Item(34, 10).Value = 'SELL'
$ws.Cells.Item(34, 11).Value = '2024-06-10 14:23:16'

# Row 35
$ws.Cells.Item(35, 8).Value = 'CANCELLED_BY_SYSTEM'
$ws.Cells.Item(35, 9).Value = 2000
$ws.Cells.Item(35, 10).Value = 'SELL'
$ws.Cells.Item(35, 11).Value = '2024-06-10 14:19:13'

# Row 36
$ws.Cells.Item(36, 8).Value = 'COMPLETED'
$ws.Cells.Item(36, 9).Value = 1000
$ws.Cells.Item(36, 10).Value = 'SELL'
$ws.Cells.Item(36, 11).Value = '2024-06-10 14:12:30'

# Row 37
$ws.Cells.Item(37, 8).Value = 'CANCELLED'
$ws.Cells.Item(37, 9).Value = 4000
$ws.Cells.Item(37, 10).Value = 'SELL'
$ws.Cells.Item(37, 11).Value = '2024-06-10 14:10:13'

# Row 38
$ws.Cells.Item(38, 8).Value = 'CANCELLED'
$ws.Cells.Item(38, 9).Value = 1500
$ws.Cells.Item(38, 10).Value = 'SELL'
$ws.Cells.Item(38, 11).Value = '2024-06-10 13:53:30'

# Row 39
$ws.Cells.Item(39, 8).Value = 'COMPLETED'
$ws.Cells.Item(39, 9).Value = 1000
$ws.Cells.Item(39, 10).Value = 'SELL'
$ws.Cells.Item(39, 11).Value = '2024-06-10 13:47:01'

# Row 40
$ws.Cells.Item(40, 8).Value = 'COMPLETED'
$ws.Cells.Item(40, 9).Value = 1000
$ws.Cells.Item(40, 10).Value = 'SELL'
$ws.Cells.Item(40, 11).Value = '2024-06-10 13:44:55'
